{"js": "// Fix the child-element order inside <w:rPr> for a handful of Pandoc\n// \"highlighting\" character styles in styles.xml. The authoring tool had\n// written <w:color> before <w:b>/<w:i>, which violates the CT_RPr sequence\n// in wml.xsd (b/bCs/i/iCs must precede color). Re-applying each style's\n// own existing bold/italic value forces the run-properties to be\n// re-serialized in schema-correct order (b, i, color) without changing any\n// actual formatting value. Only touch the toggle(s) that each style really\n// has set, so we don't materialize a new (false) <w:b>/<w:i> that wasn't\n// in the original markup.\nconst boldOnly = [\"KeywordTok\", \"ImportTok\", \"ControlFlowTok\", \"AlertTok\", \"ErrorTok\"];\nconst italicOnly = [\"CommentTok\", \"DocumentationTok\"];\nconst boldAndItalic = [\"AnnotationTok\", \"CommentVarTok\", \"InformationTok\", \"WarningTok\"];\n\nconst styles = context.document.getStyles();\n\nfunction getFont(name) {\n  const style = styles.getByNameOrNullObject(name);\n  const font = style.font;\n  font.load([\"bold\", \"italic\"]);\n  return font;\n}\n\nconst boldFonts = boldOnly.map(getFont);\nconst italicFonts = italicOnly.map(getFont);\nconst bothFonts = boldAndItalic.map(getFont);\nawait context.sync();\n\nfor (const font of boldFonts) {\n  // Re-assign the already-loaded value: no semantic change, but it makes\n  // the host rewrite <w:rPr> with elements in schema order.\n  font.bold = font.bold;\n}\nfor (const font of italicFonts) {\n  font.italic = font.italic;\n}\nfor (const font of bothFonts) {\n  font.bold = font.bold;\n  font.italic = font.italic;\n}\nawait context.sync();\n", "ps1": "# Fix the child-element order inside <w:rPr> for a handful of Pandoc\n# \"highlighting\" character styles in styles.xml. The authoring tool had\n# written <w:color> before <w:b>/<w:i>, which violates the CT_RPr sequence\n# in wml.xsd (b/bCs/i/iCs must precede color). Re-applying each style's own\n# existing bold/italic value forces the run-properties to be re-serialized\n# in schema-correct order (b, i, color) without changing any actual\n# formatting value. Only touch the toggle(s) that each style really has\n# set, so we don't materialize a new (false) <w:b>/<w:i> that wasn't in the\n# original markup.\n$d = $word.ActiveDocument\n\n$boldOnly = @(\"KeywordTok\", \"ImportTok\", \"ControlFlowTok\", \"AlertTok\", \"ErrorTok\")\n$italicOnly = @(\"CommentTok\", \"DocumentationTok\")\n$boldAndItalic = @(\"AnnotationTok\", \"CommentVarTok\", \"InformationTok\", \"WarningTok\")\n\nforeach ($name in $boldOnly) {\n    $style = $d.Styles($name)\n    $style.Font.Bold = $style.Font.Bold\n}\n\nforeach ($name in $italicOnly) {\n    $style = $d.Styles($name)\n    $style.Font.Italic = $style.Font.Italic\n}\n\nforeach ($name in $boldAndItalic) {\n    $style = $d.Styles($name)\n    $style.Font.Bold = $style.Font.Bold\n    $style.Font.Italic = $style.Font.Italic\n}\n"}
